# "Use simpler weather means" - update Climate sheet weather statistics
# and restore the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Climate")

# Updated air.temp (B), wind.2m (C) and rain.rate (D) values for each month row.
$ws.Range("B2").Value = 4.43101207056639
$ws.Range("C2").Value = 4.05891613991413
$ws.Range("D2").Value = 0.0599629009095261

$ws.Range("B3").Value = 8.23645983645984
$ws.Range("C3").Value = 3.84445591865745
$ws.Range("D3").Value = 0.0552119412831931

$ws.Range("B4").Value = 12.4492495309568
$ws.Range("C4").Value = 3.48391526295633
$ws.Range("D4").Value = 0.0702993488962998

$ws.Range("B5").Value = 16.8762259816193
$ws.Range("C5").Value = 3.15624012423227
$ws.Range("D5").Value = 0.105925308296069

$ws.Range("B6").Value = 14.4977479635841
$ws.Range("C6").Value = 3.32276959833633
$ws.Range("D6").Value = 0.128260170445409

# Restore the active cell / selection in the saved view.
$ws.Activate()
$ws.Range("F17").Select()
